$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 47) ---
$ws.Range("A47").Value = "SUBARU SUV"
$ws.Range("B47").Value = "XV_714582"
$ws.Range("C47").Value = "A"
$ws.Range("D47").Value = "AIR FILTER"
$ws.Range("E47").Value = 1000000
$ws.Range("F47").Value = 2
$ws.Range("G47").Value = "Pcs"
$ws.Range("H47").Value = "FFF"

# --- Apply the new left/right thin border to B47, then propagate the
#     resulting style to D47 and G47 via a format-only copy/paste so we
#     don't re-trigger extra intermediate style-table entries. ---
$b47 = $ws.Range("B47")
$b47.Borders(7).LineStyle = 1
$b47.Borders(7).Weight = 2
$b47.Borders(10).LineStyle = 1
$b47.Borders(10).Weight = 2

$ws.Range("B47").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("G47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Re-apply the values/text since PasteSpecial(Formats) can disturb them ---
$ws.Range("D47").Value = "AIR FILTER"
$ws.Range("G47").Value = "Pcs"

# --- Update selection to reflect the newly added row ---
$ws.Range("A46:H47").Select()
